$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 155
$ws.Range("B3").Value = 150
$ws.Range("A4").Value = "172.16.40.2"
$ws.Range("B4").Value = 21
$ws.Range("A5").Value = "1.1.1.1"
